$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated timestamps written by the pcsmote logging run (column Z = "timestamp").
$timestamps = @{
    2 = "2025-10-19T23:56:23.022296"
    3 = "2025-10-19T23:56:23.022296"
    4 = "2025-10-19T23:56:23.022851"
    5 = "2025-10-19T23:56:23.023415"
    6 = "2025-10-19T23:56:23.023931"
    7 = "2025-10-19T23:56:23.023962"
    8 = "2025-10-19T23:56:23.024511"
    9 = "2025-10-19T23:56:23.024511"
    10 = "2025-10-19T23:56:23.026570"
    11 = "2025-10-19T23:56:23.026570"
    12 = "2025-10-19T23:56:23.027572"
    13 = "2025-10-19T23:56:23.027572"
    14 = "2025-10-19T23:56:23.027572"
    15 = "2025-10-19T23:56:23.028570"
    16 = "2025-10-19T23:56:23.028570"
    17 = "2025-10-19T23:56:23.028570"
    18 = "2025-10-19T23:56:23.028570"
    19 = "2025-10-19T23:56:23.029572"
    20 = "2025-10-19T23:56:23.029572"
    21 = "2025-10-19T23:56:23.029572"
    22 = "2025-10-19T23:56:23.030573"
    23 = "2025-10-19T23:56:23.030573"
    24 = "2025-10-19T23:56:23.030573"
    25 = "2025-10-19T23:56:23.030573"
    26 = "2025-10-19T23:56:23.030573"
    27 = "2025-10-19T23:56:23.031572"
    28 = "2025-10-19T23:56:23.031572"
    29 = "2025-10-19T23:56:23.031572"
    30 = "2025-10-19T23:56:23.031572"
    31 = "2025-10-19T23:56:23.032572"
    32 = "2025-10-19T23:56:23.032572"
    33 = "2025-10-19T23:56:23.032572"
    34 = "2025-10-19T23:56:23.032572"
    35 = "2025-10-19T23:56:23.033572"
    36 = "2025-10-19T23:56:23.033572"
    37 = "2025-10-19T23:56:23.033572"
    38 = "2025-10-19T23:56:23.033572"
    39 = "2025-10-19T23:56:23.034574"
    40 = "2025-10-19T23:56:23.034574"
    41 = "2025-10-19T23:56:23.034574"
    42 = "2025-10-19T23:56:23.034574"
    43 = "2025-10-19T23:56:23.035574"
    44 = "2025-10-19T23:56:23.035574"
    45 = "2025-10-19T23:56:23.035574"
    46 = "2025-10-19T23:56:23.110519"
    47 = "2025-10-19T23:56:23.111519"
    48 = "2025-10-19T23:56:23.111519"
    49 = "2025-10-19T23:56:23.111519"
    50 = "2025-10-19T23:56:23.111519"
    51 = "2025-10-19T23:56:23.112518"
    52 = "2025-10-19T23:56:23.112518"
    53 = "2025-10-19T23:56:23.112518"
    54 = "2025-10-19T23:56:23.112518"
    55 = "2025-10-19T23:56:23.112518"
    56 = "2025-10-19T23:56:23.113518"
    57 = "2025-10-19T23:56:23.113518"
    58 = "2025-10-19T23:56:23.113518"
    59 = "2025-10-19T23:56:23.113518"
    60 = "2025-10-19T23:56:23.113518"
    61 = "2025-10-19T23:56:23.141252"
    62 = "2025-10-19T23:56:23.141777"
    63 = "2025-10-19T23:56:23.142775"
    64 = "2025-10-19T23:56:23.142775"
    65 = "2025-10-19T23:56:23.143792"
    66 = "2025-10-19T23:56:23.144810"
    67 = "2025-10-19T23:56:23.144810"
    68 = "2025-10-19T23:56:23.145780"
    69 = "2025-10-19T23:56:23.145780"
    70 = "2025-10-19T23:56:23.145780"
    71 = "2025-10-19T23:56:23.145780"
    72 = "2025-10-19T23:56:23.145780"
    73 = "2025-10-19T23:56:23.146776"
    74 = "2025-10-19T23:56:23.146776"
    75 = "2025-10-19T23:56:23.227981"
    76 = "2025-10-19T23:56:23.228979"
    77 = "2025-10-19T23:56:23.228979"
    78 = "2025-10-19T23:56:23.228979"
    79 = "2025-10-19T23:56:23.228979"
    80 = "2025-10-19T23:56:23.228979"
    81 = "2025-10-19T23:56:23.229978"
    82 = "2025-10-19T23:56:23.230980"
    83 = "2025-10-19T23:56:23.230980"
    84 = "2025-10-19T23:56:23.230980"
    85 = "2025-10-19T23:56:23.230980"
    86 = "2025-10-19T23:56:23.230980"
    87 = "2025-10-19T23:56:23.231978"
    88 = "2025-10-19T23:56:23.231978"
    89 = "2025-10-19T23:56:23.232979"
    90 = "2025-10-19T23:56:23.232979"
    91 = "2025-10-19T23:56:23.232979"
    92 = "2025-10-19T23:56:23.232979"
    93 = "2025-10-19T23:56:23.233979"
    94 = "2025-10-19T23:56:23.233979"
    95 = "2025-10-19T23:56:23.233979"
    96 = "2025-10-19T23:56:23.233979"
    97 = "2025-10-19T23:56:23.233979"
    98 = "2025-10-19T23:56:23.234978"
    99 = "2025-10-19T23:56:23.234978"
    100 = "2025-10-19T23:56:23.234978"
    101 = "2025-10-19T23:56:23.234978"
    102 = "2025-10-19T23:56:23.234978"
    103 = "2025-10-19T23:56:23.262050"
    104 = "2025-10-19T23:56:23.262050"
    105 = "2025-10-19T23:56:23.262583"
    106 = "2025-10-19T23:56:23.262583"
    107 = "2025-10-19T23:56:23.263158"
    108 = "2025-10-19T23:56:23.263158"
    109 = "2025-10-19T23:56:23.263707"
    110 = "2025-10-19T23:56:23.264693"
    111 = "2025-10-19T23:56:23.264693"
    112 = "2025-10-19T23:56:23.264693"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item([int]$row, 26).Value = $timestamps[$row]
}

